# Take Home Ex 4 - apply edits to Steps.pptx
# EMU -> points helper (PowerPoint Shape geometry properties are expressed in points)
function EMU([double]$emu) { return $emu / 12700.0 }

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1: move the screenshot picture up/left a bit and add a new red
# highlight rectangle (solid line) next to it.
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$pic1 = $s1.Shapes.Item("Picture 3")
$pic1.Left = EMU(1093579)
$pic1.Top  = EMU(1090157)

$rect1 = $s1.Shapes.AddShape(1, (EMU 1110831), (EMU 2829464), (EMU 398792), (EMU 3260785))
$rect1.Name = "Rectangle 4"
$rect1.Fill.Visible = 0
$rect1.Line.Visible = -1
$rect1.Line.ForeColor.RGB = 255
$rect1.Line.Weight = 3
$rect1.Line.DashStyle = 1
$rect1.TextFrame.VerticalAnchor = 3
$rect1.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# ---------------------------------------------------------------------------
# Slide 3: remove the first red highlight rectangle and change the
# remaining one's line from dashed to solid (with a small resize).
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item("Rectangle 3").Delete()
$rect3b = $s3.Shapes.Item("Rectangle 4")
$rect3b.Left   = EMU(1217603)
$rect3b.Top    = EMU(3200589)
$rect3b.Width  = EMU(1226191)
$rect3b.Height = EMU(551902)
$rect3b.Line.DashStyle = 1

# ---------------------------------------------------------------------------
# Slide 4: drop the two old highlight rectangles, move the caption textbox
# up in the stacking order, and draw two new solid-line highlight
# rectangles in adjusted positions.
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item("Rectangle 3").Delete()
$s4.Shapes.Item("Rectangle 4").Delete()

$tb4 = $s4.Shapes.Item("TextBox 5")
$tb4.ZOrder(3)  # msoSendBackward: move it just after the picture

$rect4a = $s4.Shapes.AddShape(1, (EMU 4209690), (EMU 3200585), (EMU 698739), (EMU 185279))
$rect4a.Name = "Rectangle 6"
$rect4a.Fill.Visible = 0
$rect4a.Line.Visible = -1
$rect4a.Line.ForeColor.RGB = 255
$rect4a.Line.Weight = 3
$rect4a.Line.DashStyle = 1
$rect4a.TextFrame.VerticalAnchor = 3
$rect4a.TextFrame.TextRange.ParagraphFormat.Alignment = 2

$rect4b = $s4.Shapes.AddShape(1, (EMU 4090199), (EMU 2533511), (EMU 205755), (EMU 185279))
$rect4b.Name = "Rectangle 7"
$rect4b.Fill.Visible = 0
$rect4b.Line.Visible = -1
$rect4b.Line.ForeColor.RGB = 255
$rect4b.Line.Weight = 3
$rect4b.Line.DashStyle = 1
$rect4b.TextFrame.VerticalAnchor = 3
$rect4b.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# ---------------------------------------------------------------------------
# Slide 9: remove the duplicate caption textbox.
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$s9.Shapes.Item("TextBox 4").Delete()

# ---------------------------------------------------------------------------
# Slide master / layouts: the cached "today" date placeholder text bumped
# from 12/2/2022 to 13/2/2022.
# ---------------------------------------------------------------------------
$master = $p.SlideMaster
$masterDt = $master.HeadersFooters.DateAndTime
$masterDt.Text = "13/2/2022"
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    $layoutDt = $layout.HeadersFooters.DateAndTime
    $layoutDt.Text = "13/2/2022"
}
